# Update "想去人数" (F column) / "最低票价" (G column) values in the
# "展览" and "全部类型" worksheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> @{ Column = NewValue }
$updates = @{
    2  = @{ G = 65 }
    5  = @{ F = 2809 }
    9  = @{ F = 1506 }
    10 = @{ F = 35 }
    13 = @{ F = 1254 }
    15 = @{ F = 388 }
    16 = @{ F = 333 }
    18 = @{ F = 45 }
    19 = @{ F = 112 }
    22 = @{ F = 2776 }
    25 = @{ F = 50 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $cellAddr = "$col$row"
            $ws.Range($cellAddr).Value = $cols[$col]
        }
    }
}
